$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx")

if ($found) {
    $jupParagraph = $rng.Paragraphs(1)

    # The paragraph right before it is the blank paragraph that should be
    # removed along with it; the paragraph right after it is the copyright
    # notice paragraph, which should also be removed.
    $prevParagraph = $jupParagraph.Previous()
    $nextParagraph = $jupParagraph.Next()

    $delStart = $prevParagraph.Range.Start
    $delEnd = $nextParagraph.Range.End

    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
